# Apply scheduled-runner Sheets update: refresh profit-calc cells (H/I/J/K/L/M/N)
# across ALC, ARM, CRP, CUL, GSM, LTW, WVR per latest market-price snapshot.
$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
# Row 3
$wsALC.Range("H3").Value = 29578.5
$wsALC.Range("J3").Value = 29578.5
$wsALC.Range("L3").Value = 29578.5
$wsALC.Range("N3").Value = -29806.5

# Row 11
$wsALC.Range("H11").Value = 550095.4
$wsALC.Range("I11").Value = 550095.4
$wsALC.Range("K11").Value = 550095.4
$wsALC.Range("M11").Value = -549955.4

# Row 102
$wsALC.Range("H102").Value = 29578.5
$wsALC.Range("J102").Value = 29578.5
$wsALC.Range("L102").Value = 29578.5
$wsALC.Range("N102").Value = -36068.5

# Row 118
$wsALC.Range("H118").Value = 11667.667
$wsALC.Range("I118").Value = 14651.286
$wsALC.Range("J118").Value = 1225
$wsALC.Range("K118").Value = 43953.858
$wsALC.Range("L118").Value = 3675
$wsALC.Range("M118").Value = -42296.858
$wsALC.Range("N118").Value = -6989

# Row 132
$wsALC.Range("H132").Value = 5819203
$wsALC.Range("I132").Value = 6415695.5
$wsALC.Range("J132").Value = 3402.75
$wsALC.Range("K132").Value = 19247086.5
$wsALC.Range("L132").Value = 10208.25
$wsALC.Range("M132").Value = -19244556.5
$wsALC.Range("N132").Value = -15268.25

$wsARM = $wb.Worksheets.Item("ARM")
# Row 26
$wsARM.Range("H26").Value = 5421.1113
$wsARM.Range("I26").Value = 5197.5
$wsARM.Range("J26").Value = 5600
$wsARM.Range("K26").Value = 5197.5
$wsARM.Range("L26").Value = 5600
$wsARM.Range("M26").Value = -4867.5
$wsARM.Range("N26").Value = -6260

# Row 32
$wsARM.Range("H32").Value = 24140.389
$wsARM.Range("I32").Value = 4435.894
$wsARM.Range("J32").Value = 142367.36
$wsARM.Range("K32").Value = 4435.894
$wsARM.Range("L32").Value = 142367.36
$wsARM.Range("M32").Value = -4148.894
$wsARM.Range("N32").Value = -142941.36

# Row 45
$wsARM.Range("H45").Value = 53621.5
$wsARM.Range("I45").Value = 69901.47
$wsARM.Range("J45").Value = 4781.6
$wsARM.Range("K45").Value = 69901.47
$wsARM.Range("L45").Value = 4781.6
$wsARM.Range("M45").Value = -69524.47
$wsARM.Range("N45").Value = -5535.6

# Row 46
$wsARM.Range("H46").Value = 3664.5715
$wsARM.Range("I46").Value = 1750
$wsARM.Range("J46").Value = 4430.4
$wsARM.Range("K46").Value = 1750
$wsARM.Range("L46").Value = 4430.4
$wsARM.Range("M46").Value = -1431
$wsARM.Range("N46").Value = -5068.4

# Row 74
$wsARM.Range("H74").Value = 1373.9333
$wsARM.Range("I74").Value = 661.65216
$wsARM.Range("K74").Value = 661.65216
$wsARM.Range("M74").Value = 212.34784

# Row 77
$wsARM.Range("H77").Value = 1373.9333
$wsARM.Range("I77").Value = 661.65216
$wsARM.Range("K77").Value = 3308.2608
$wsARM.Range("M77").Value = 1059.7392

# Row 110
$wsARM.Range("H110").Value = 62626570
$wsARM.Range("I110").Value = 71572936
$wsARM.Range("J110").Value = 1999.5
$wsARM.Range("K110").Value = 71572936
$wsARM.Range("L110").Value = 1999.5
$wsARM.Range("M110").Value = -71570891
$wsARM.Range("N110").Value = -6089.5

# Row 132
$wsARM.Range("H132").Value = 1647.3125
$wsARM.Range("I132").Value = 1058.9697
$wsARM.Range("J132").Value = 2941.6667
$wsARM.Range("K132").Value = 3176.9091
$wsARM.Range("L132").Value = 8825.000100000001
$wsARM.Range("M132").Value = -646.9091000000003
$wsARM.Range("N132").Value = -13885.0001

$wsCRP = $wb.Worksheets.Item("CRP")
# Row 16
$wsCRP.Range("H16").Value = 1159
$wsCRP.Range("J16").Value = 995
$wsCRP.Range("L16").Value = 995
$wsCRP.Range("N16").Value = -1569

# Row 35
$wsCRP.Range("H35").Value = 2472.2727
$wsCRP.Range("I35").Value = 1561.875
$wsCRP.Range("J35").Value = 4900
$wsCRP.Range("K35").Value = 1561.875
$wsCRP.Range("L35").Value = 4900
$wsCRP.Range("M35").Value = -1267.875
$wsCRP.Range("N35").Value = -5488

# Row 74
$wsCRP.Range("H74").Value = 23138.5
$wsCRP.Range("J74").Value = 23138.5
$wsCRP.Range("L74").Value = 23138.5
$wsCRP.Range("N74").Value = -24886.5

# Row 77
$wsCRP.Range("H77").Value = 23138.5
$wsCRP.Range("J77").Value = 23138.5
$wsCRP.Range("L77").Value = 69415.5
$wsCRP.Range("N77").Value = -78151.5

# Row 105
$wsCRP.Range("H105").Value = 2284.5293
$wsCRP.Range("I105").Value = 2366.5386
$wsCRP.Range("K105").Value = 2366.5386
$wsCRP.Range("M105").Value = -619.5385999999999

# Row 112
$wsCRP.Range("H112").Value = 0
$wsCRP.Range("J112").Value = 0
$wsCRP.Range("L112").Value = 0
$wsCRP.Range("N112").ClearContents()

# Row 113
$wsCRP.Range("H113").Value = 1159
$wsCRP.Range("J113").Value = 995
$wsCRP.Range("L113").Value = 995
$wsCRP.Range("N113").Value = -5335

# Row 132
$wsCRP.Range("H132").Value = 3805.9736
$wsCRP.Range("I132").Value = 4039.1365
$wsCRP.Range("J132").Value = 3485.375
$wsCRP.Range("K132").Value = 12117.4095
$wsCRP.Range("L132").Value = 10456.125
$wsCRP.Range("M132").Value = -9587.4095
$wsCRP.Range("N132").Value = -15516.125

# Row 134
$wsCRP.Range("H134").Value = 1238.1945
$wsCRP.Range("I134").Value = 1157.8148
$wsCRP.Range("K134").Value = 3473.4444
$wsCRP.Range("M134").Value = -938.4444000000003

$wsCUL = $wb.Worksheets.Item("CUL")
# Row 5
$wsCUL.Range("H5").Value = 1427.7368
$wsCUL.Range("I5").Value = 759.55
$wsCUL.Range("J5").Value = 2170.1667
$wsCUL.Range("K5").Value = 2278.65
$wsCUL.Range("L5").Value = 6510.500100000001
$wsCUL.Range("M5").Value = -2166.65
$wsCUL.Range("N5").Value = -6734.500100000001

# Row 96
$wsCUL.Range("H96").Value = 3866.6667
$wsCUL.Range("J96").Value = 3866.6667
$wsCUL.Range("L96").Value = 11600.0001
$wsCUL.Range("N96").Value = -15718.0001

# Row 105
$wsCUL.Range("H105").Value = 9811.6
$wsCUL.Range("I105").Value = 3000
$wsCUL.Range("J105").Value = 11514.5
$wsCUL.Range("K105").Value = 9000
$wsCUL.Range("L105").Value = 34543.5
$wsCUL.Range("M105").Value = -6379
$wsCUL.Range("N105").Value = -39785.5

# Row 131
$wsCUL.Range("H131").Value = 808.41
$wsCUL.Range("J131").Value = 856.52747
$wsCUL.Range("L131").Value = 2569.58241
$wsCUL.Range("N131").Value = -12649.58241

# Row 132
$wsCUL.Range("H132").Value = 1496.7646
$wsCUL.Range("I132").Value = 912
$wsCUL.Range("J132").Value = 1740.4166
$wsCUL.Range("K132").Value = 8208
$wsCUL.Range("L132").Value = 15663.7494
$wsCUL.Range("M132").Value = -5678
$wsCUL.Range("N132").Value = -20723.7494

# Row 135
$wsCUL.Range("H135").Value = 1427.7368
$wsCUL.Range("I135").Value = 759.55
$wsCUL.Range("J135").Value = 2170.1667
$wsCUL.Range("K135").Value = 6835.95
$wsCUL.Range("L135").Value = 19531.5003
$wsCUL.Range("M135").Value = -4300.95
$wsCUL.Range("N135").Value = -24601.5003

$wsGSM = $wb.Worksheets.Item("GSM")
# Row 70
$wsGSM.Range("H70").Value = 76793.53999999999
$wsGSM.Range("I70").Value = 158192.47
$wsGSM.Range("J70").Value = 6247.8
$wsGSM.Range("K70").Value = 158192.47
$wsGSM.Range("L70").Value = 6247.8
$wsGSM.Range("M70").Value = -157922.47
$wsGSM.Range("N70").Value = -6787.8

# Row 73
$wsGSM.Range("H73").Value = 76793.53999999999
$wsGSM.Range("I73").Value = 158192.47
$wsGSM.Range("J73").Value = 6247.8
$wsGSM.Range("K73").Value = 158192.47
$wsGSM.Range("L73").Value = 6247.8
$wsGSM.Range("M73").Value = -157256.47
$wsGSM.Range("N73").Value = -8119.8

$wsLTW = $wb.Worksheets.Item("LTW")
# Row 115
$wsLTW.Range("H115").Value = 0
$wsLTW.Range("J115").Value = 0
$wsLTW.Range("L115").Value = 0
$wsLTW.Range("N115").ClearContents()

$wsWVR = $wb.Worksheets.Item("WVR")
# Row 75
$wsWVR.Range("H75").Value = 10000000
$wsWVR.Range("J75").Value = 10000000
$wsWVR.Range("L75").Value = 10000000
$wsWVR.Range("N75").Value = -10001872

# Row 78
$wsWVR.Range("H78").Value = 10000000
$wsWVR.Range("J78").Value = 10000000
$wsWVR.Range("L78").Value = 30000000
$wsWVR.Range("N78").Value = -30009360

# Row 82
$wsWVR.Range("H82").Value = 16666.666
$wsWVR.Range("J82").Value = 16666.666
$wsWVR.Range("L82").Value = 16666.666
$wsWVR.Range("N82").Value = -17432.666

# Row 85
$wsWVR.Range("H85").Value = 16666.666
$wsWVR.Range("J85").Value = 16666.666
$wsWVR.Range("L85").Value = 16666.666
$wsWVR.Range("N85").Value = -19318.666

# Row 109
$wsWVR.Range("H109").Value = 32500
$wsWVR.Range("J109").Value = 32500
$wsWVR.Range("L109").Value = 32500
$wsWVR.Range("N109").Value = -35274

# Row 116
$wsWVR.Range("H116").Value = 49970
$wsWVR.Range("J116").Value = 49970
$wsWVR.Range("L116").Value = 49970
$wsWVR.Range("N116").Value = -59148

# Row 132
$wsWVR.Range("H132").Value = 2201.742
$wsWVR.Range("I132").Value = 1766.0204
$wsWVR.Range("J132").Value = 3844.077
$wsWVR.Range("K132").Value = 5298.0612
$wsWVR.Range("L132").Value = 11532.231
$wsWVR.Range("M132").Value = -2768.0612
$wsWVR.Range("N132").Value = -16592.231
